$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 500
$ws.Range("J2").Value = 500
$ws.Range("L2").Value = 500
$ws.Range("N2").Value = -726
$ws.Range("H4").Value = 63.4
$ws.Range("J4").Value = 40
$ws.Range("L4").Value = 40
$ws.Range("N4").Value = -268
$ws.Range("H5").Value = 76
$ws.Range("I5").Value = 81.333336
$ws.Range("J5").Value = 60
$ws.Range("K5").Value = 81.333336
$ws.Range("L5").Value = 60
$ws.Range("M5").Value = 33.666664
$ws.Range("N5").Value = -290
$ws.Range("H9").Value = 99.583336
$ws.Range("I9").Value = 110
$ws.Range("J9").Value = 47.5
$ws.Range("K9").Value = 110
$ws.Range("L9").Value = 47.5
$ws.Range("M9").Value = 59
$ws.Range("N9").Value = -385.5
$ws.Range("H28").Value = 2481.5
$ws.Range("I28").Value = 2477.9
$ws.Range("K28").Value = 2477.9
$ws.Range("M28").Value = -1992.9
$ws.Range("H32").Value = 14291586
$ws.Range("J32").Value = 14291586
$ws.Range("L32").Value = 14291586
$ws.Range("N32").Value = -14292238
$ws.Range("H40").Value = 4898.9
$ws.Range("I40").Value = 1459
$ws.Range("J40").Value = 5758.875
$ws.Range("K40").Value = 1459
$ws.Range("L40").Value = 5758.875
$ws.Range("M40").Value = -1284
$ws.Range("N40").Value = -6108.875
$ws.Range("H41").Value = 1020.2381
$ws.Range("I41").Value = 692.1818
$ws.Range("J41").Value = 1381.1
$ws.Range("K41").Value = 692.1818
$ws.Range("L41").Value = 1381.1
$ws.Range("M41").Value = -252.1818
$ws.Range("N41").Value = -2261.1
$ws.Range("H48").Value = 4499.5
$ws.Range("J48").Value = 4499.5
$ws.Range("L48").Value = 13498.5
$ws.Range("N48").Value = -14082.5
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("H55").Value = 630.05554
$ws.Range("I55").Value = 635.875
$ws.Range("K55").Value = 635.875
$ws.Range("M55").Value = -421.875
$ws.Range("H56").Value = 4499.5
$ws.Range("J56").Value = 4499.5
$ws.Range("L56").Value = 13498.5
$ws.Range("N56").Value = -14566.5
$ws.Range("H62").Value = 4638.4
$ws.Range("J62").Value = 4731.3335
$ws.Range("L62").Value = 4731.3335
$ws.Range("N62").Value = -5979.3335
$ws.Range("H65").Value = 4638.4
$ws.Range("J65").Value = 4731.3335
$ws.Range("L65").Value = 23656.6675
$ws.Range("N65").Value = -29896.6675
$ws.Range("H69").Value = 28830.295
$ws.Range("J69").Value = 41303.637
$ws.Range("L69").Value = 123910.911
$ws.Range("N69").Value = -125658.911
$ws.Range("H72").Value = 28830.295
$ws.Range("J72").Value = 41303.637
$ws.Range("L72").Value = 371732.733
$ws.Range("N72").Value = -380468.733
$ws.Range("H87").Value = 19999.76
$ws.Range("J87").Value = 19999.76
$ws.Range("L87").Value = 19999.76
$ws.Range("N87").Value = -22495.76
$ws.Range("H88").Value = 1563.2354
$ws.Range("J88").Value = 656
$ws.Range("L88").Value = 656
$ws.Range("N88").Value = -1468
$ws.Range("H90").Value = 19999.76
$ws.Range("J90").Value = 19999.76
$ws.Range("L90").Value = 59999.28
$ws.Range("N90").Value = -72479.28
$ws.Range("H91").Value = 1563.2354
$ws.Range("J91").Value = 656
$ws.Range("L91").Value = 656
$ws.Range("N91").Value = -3464
$ws.Range("H92").Value = 1107.6086
$ws.Range("J92").Value = 1157.5
$ws.Range("L92").Value = 1157.5
$ws.Range("N92").Value = -3653.5
$ws.Range("H98").Value = 3500
$ws.Range("I98").Value = 2000
$ws.Range("K98").Value = 2000
$ws.Range("M98").Value = -502
$ws.Range("H99").Value = 66777.69
$ws.Range("I99").Value = 325.77777
$ws.Range("K99").Value = 977.33331
$ws.Range("M99").Value = 520.66669
$ws.Range("H100").Value = 5871.24
$ws.Range("I100").Value = 4664.6665
$ws.Range("J100").Value = 6985
$ws.Range("K100").Value = 4664.6665
$ws.Range("L100").Value = 6985
$ws.Range("M100").Value = -4123.6665
$ws.Range("N100").Value = -8067
$ws.Range("H112").Value = 1465.3934
$ws.Range("J112").Value = 1478.138
$ws.Range("L112").Value = 4434.414
$ws.Range("N112").Value = -6650.414
$ws.Range("H122").Value = 3500
$ws.Range("I122").Value = 2000
$ws.Range("K122").Value = 6000
$ws.Range("M122").Value = -3550
$ws.Range("H132").Value = 18099.838
$ws.Range("I132").Value = 1545.6666
$ws.Range("K132").Value = 4636.9998
$ws.Range("M132").Value = -2106.9998
$ws.Range("H135").Value = 3241
$ws.Range("I135").Value = 2846.8572
$ws.Range("K135").Value = 25621.7148
$ws.Range("M135").Value = -23086.7148
$ws.Range("H138").Value = 3171.2307
$ws.Range("I138").Value = 2075.75
$ws.Range("J138").Value = 4110.2144
$ws.Range("K138").Value = 6227.25
$ws.Range("L138").Value = 12330.6432
$ws.Range("M138").Value = -1087.25
$ws.Range("N138").Value = -22610.6432
$ws.Range("H141").Value = 3524.8
$ws.Range("I141").Value = 3300.2
$ws.Range("J141").Value = 4198.6
$ws.Range("K141").Value = 9900.599999999999
$ws.Range("L141").Value = 12595.8
$ws.Range("M141").Value = -4720.599999999999
$ws.Range("N141").Value = -22955.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 21278.666
$ws.Range("I2").Value = 34085.89
$ws.Range("J2").Value = 2067.8333
$ws.Range("K2").Value = 34085.89
$ws.Range("L2").Value = 2067.8333
$ws.Range("M2").Value = -33972.89
$ws.Range("N2").Value = -2293.8333
$ws.Range("H32").Value = 9939.091
$ws.Range("I32").Value = 7181.525
$ws.Range("K32").Value = 7181.525
$ws.Range("M32").Value = -6894.525
$ws.Range("H45").Value = 2821.7273
$ws.Range("J45").Value = 5004.3335
$ws.Range("L45").Value = 5004.3335
$ws.Range("N45").Value = -5758.3335
$ws.Range("H102").Value = 2200.0527
$ws.Range("I102").Value = 2046.5883
$ws.Range("J102").Value = 3504.5
$ws.Range("K102").Value = 2046.5883
$ws.Range("L102").Value = 3504.5
$ws.Range("M102").Value = -424.5882999999999
$ws.Range("N102").Value = -6748.5
$ws.Range("H116").Value = 21278.666
$ws.Range("I116").Value = 34085.89
$ws.Range("J116").Value = 2067.8333
$ws.Range("K116").Value = 34085.89
$ws.Range("L116").Value = 2067.8333
$ws.Range("M116").Value = -31791.89
$ws.Range("N116").Value = -6655.8333
$ws.Range("H122").Value = 4537.44
$ws.Range("I122").Value = 2760.9412
$ws.Range("K122").Value = 8282.8236
$ws.Range("M122").Value = -5832.8236
$ws.Range("H124").Value = 29996
$ws.Range("J124").Value = 29996
$ws.Range("L124").Value = 29996
$ws.Range("N124").Value = -39816
$ws.Range("H132").Value = 2568.8684
$ws.Range("I132").Value = 1871.0358
$ws.Range("J132").Value = 4522.8
$ws.Range("K132").Value = 5613.107400000001
$ws.Range("L132").Value = 13568.4
$ws.Range("M132").Value = -3083.107400000001
$ws.Range("N132").Value = -18628.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 21278.666
$ws.Range("I3").Value = 34085.89
$ws.Range("J3").Value = 2067.8333
$ws.Range("K3").Value = 34085.89
$ws.Range("L3").Value = 2067.8333
$ws.Range("M3").Value = -33971.89
$ws.Range("N3").Value = -2295.8333
$ws.Range("H22").Value = 762.25
$ws.Range("I22").Value = 659.8
$ws.Range("J22").Value = 933
$ws.Range("K22").Value = 659.8
$ws.Range("L22").Value = 933
$ws.Range("M22").Value = -486.8
$ws.Range("N22").Value = -1279
$ws.Range("H57").Value = 99995.336
$ws.Range("J57").Value = 99995.5
$ws.Range("L57").Value = 99995.5
$ws.Range("N57").Value = -101435.5
$ws.Range("H58").Value = 41434.5
$ws.Range("I58").Value = 44869
$ws.Range("K58").Value = 44869
$ws.Range("M58").Value = -44575
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("H60").Value = 72498.75
$ws.Range("I60").Value = 59995
$ws.Range("K60").Value = 59995
$ws.Range("M60").Value = -59396
$ws.Range("H76").Value = 31999.8
$ws.Range("I76").Value = 33999.668
$ws.Range("K76").Value = 33999.668
$ws.Range("M76").Value = -33684.668
$ws.Range("H79").Value = 31999.8
$ws.Range("I79").Value = 33999.668
$ws.Range("K79").Value = 33999.668
$ws.Range("M79").Value = -32907.668
$ws.Range("H86").Value = 2629.862
$ws.Range("I86").Value = 1197.8636
$ws.Range("K86").Value = 1197.8636
$ws.Range("M86").Value = -74.86359999999991
$ws.Range("H89").Value = 2629.862
$ws.Range("I89").Value = 1197.8636
$ws.Range("K89").Value = 5989.317999999999
$ws.Range("M89").Value = -373.3179999999993
$ws.Range("H94").Value = 3052
$ws.Range("I94").Value = 2941.7896
$ws.Range("J94").Value = 3401
$ws.Range("K94").Value = 2941.7896
$ws.Range("L94").Value = 3401
$ws.Range("M94").Value = -2490.7896
$ws.Range("N94").Value = -4303
$ws.Range("H107").Value = 1771.6154
$ws.Range("I107").Value = 1592.5454
$ws.Range("J107").Value = 2756.5
$ws.Range("K107").Value = 1592.5454
$ws.Range("L107").Value = 2756.5
$ws.Range("M107").Value = 327.4546
$ws.Range("N107").Value = -6596.5
$ws.Range("H132").Value = 70110.94500000001
$ws.Range("J132").Value = 70110.94500000001
$ws.Range("L132").Value = 70110.94500000001
$ws.Range("N132").Value = -80230.94500000001
$ws.Range("H133").Value = 55000
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H134").Value = 4043.3333
$ws.Range("I134").Value = 3144.0908
$ws.Range("K134").Value = 9432.2724
$ws.Range("M134").Value = -6897.2724
$ws.Range("H136").Value = 99995.336
$ws.Range("J136").Value = 99995.5
$ws.Range("L136").Value = 99995.5
$ws.Range("N136").Value = -110195.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 187.05556
$ws.Range("I7").Value = 27
$ws.Range("J7").Value = 248.61539
$ws.Range("K7").Value = 27
$ws.Range("L7").Value = 248.61539
$ws.Range("M7").Value = 86
$ws.Range("N7").Value = -474.61539
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H52").Value = 65702
$ws.Range("J52").Value = 66842.39999999999
$ws.Range("L52").Value = 66842.39999999999
$ws.Range("N52").Value = -67430.39999999999
$ws.Range("H62").Value = 5616.8184
$ws.Range("I62").Value = 5678.7
$ws.Range("K62").Value = 5678.7
$ws.Range("M62").Value = -5054.7
$ws.Range("H65").Value = 5616.8184
$ws.Range("I65").Value = 5678.7
$ws.Range("K65").Value = 28393.5
$ws.Range("M65").Value = -25273.5
$ws.Range("H107").Value = 11158.6
$ws.Range("I107").Value = 1600
$ws.Range("J107").Value = 17531
$ws.Range("K107").Value = 1600
$ws.Range("L107").Value = 17531
$ws.Range("M107").Value = 320
$ws.Range("N107").Value = -21371
$ws.Range("H132").Value = 2442.2778
$ws.Range("I132").Value = 2442.2778
$ws.Range("K132").Value = 7326.8334
$ws.Range("M132").Value = -4796.8334
$ws.Range("H134").Value = 3800.6938
$ws.Range("I134").Value = 2540.1516
$ws.Range("K134").Value = 7620.4548
$ws.Range("M134").Value = -5085.4548
$ws.Range("H139").Value = 89992.14
$ws.Range("J139").Value = 89992.14
$ws.Range("L139").Value = 89992.14
$ws.Range("N139").Value = -100272.14
$ws.Range("H141").Value = 119000
$ws.Range("J141").Value = 119000
$ws.Range("L141").Value = 119000
$ws.Range("N141").Value = -129360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 204.22223
$ws.Range("I6").Value = 119.85714
$ws.Range("J6").Value = 499.5
$ws.Range("K6").Value = 359.57142
$ws.Range("L6").Value = 1498.5
$ws.Range("M6").Value = -246.57142
$ws.Range("N6").Value = -1724.5
$ws.Range("H37").Value = 199994.5
$ws.Range("J37").Value = 199994.5
$ws.Range("L37").Value = 599983.5
$ws.Range("N37").Value = -600207.5
$ws.Range("H40").Value = 57.555557
$ws.Range("I40").Value = 56.666668
$ws.Range("J40").Value = 59.333332
$ws.Range("K40").Value = 226.666672
$ws.Range("L40").Value = 237.333328
$ws.Range("M40").Value = -157.666672
$ws.Range("N40").Value = -375.333328
$ws.Range("H51").Value = 2517.25
$ws.Range("I51").Value = 2517.25
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 7551.75
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -7091.75
$ws.Range("N51").ClearContents()
$ws.Range("H81").Value = 2626.25
$ws.Range("J81").Value = 2626.25
$ws.Range("L81").Value = 7878.75
$ws.Range("N81").Value = -10124.75
$ws.Range("H84").Value = 2626.25
$ws.Range("J84").Value = 2626.25
$ws.Range("L84").Value = 23636.25
$ws.Range("N84").Value = -34868.25
$ws.Range("H107").Value = 441.8
$ws.Range("J107").Value = 502.5
$ws.Range("L107").Value = 1507.5
$ws.Range("N107").Value = -5347.5
$ws.Range("H113").Value = 6287.737
$ws.Range("I113").Value = 9894.546
$ws.Range("K113").Value = 29683.638
$ws.Range("M113").Value = -27513.638
$ws.Range("H131").Value = 8180.0586
$ws.Range("I131").Value = 4754.8335
$ws.Range("J131").Value = 10048.363
$ws.Range("K131").Value = 14264.5005
$ws.Range("L131").Value = 30145.089
$ws.Range("M131").Value = -9224.500499999998
$ws.Range("N131").Value = -40225.089
$ws.Range("H132").Value = 1236.75
$ws.Range("I132").Value = 1183.3334
$ws.Range("J132").Value = 1397
$ws.Range("K132").Value = 10650.0006
$ws.Range("L132").Value = 12573
$ws.Range("M132").Value = -8120.000599999999
$ws.Range("N132").Value = -17633
$ws.Range("H138").Value = 7732.4
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 7732.4
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 23197.2
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -33477.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 4000254.5
$ws.Range("I2").Value = 241
$ws.Range("J2").Value = 8333602.5
$ws.Range("K2").Value = 241
$ws.Range("L2").Value = 8333602.5
$ws.Range("M2").Value = -128
$ws.Range("N2").Value = -8333828.5
$ws.Range("H33").Value = 9028.833000000001
$ws.Range("J33").Value = 9028.833000000001
$ws.Range("L33").Value = 9028.833000000001
$ws.Range("N33").Value = -9532.833000000001
$ws.Range("H40").Value = 6999.75
$ws.Range("I40").Value = 4000
$ws.Range("J40").Value = 7999.6665
$ws.Range("K40").Value = 4000
$ws.Range("L40").Value = 7999.6665
$ws.Range("M40").Value = -3849
$ws.Range("N40").Value = -8301.666499999999
$ws.Range("H44").Value = 20552.25
$ws.Range("J44").Value = 24436.334
$ws.Range("L44").Value = 24436.334
$ws.Range("N44").Value = -25628.334
$ws.Range("H80").Value = 161235.72
$ws.Range("I80").Value = 279452.75
$ws.Range("J80").Value = 3613
$ws.Range("K80").Value = 279452.75
$ws.Range("L80").Value = 3613
$ws.Range("M80").Value = -278454.75
$ws.Range("N80").Value = -5609
$ws.Range("H83").Value = 161235.72
$ws.Range("I83").Value = 279452.75
$ws.Range("J83").Value = 3613
$ws.Range("K83").Value = 1397263.75
$ws.Range("L83").Value = 18065
$ws.Range("M83").Value = -1392271.75
$ws.Range("N83").Value = -28049
$ws.Range("H107").Value = 665
$ws.Range("I107").Value = 384.33334
$ws.Range("J107").Value = 1001.8
$ws.Range("K107").Value = 384.33334
$ws.Range("L107").Value = 1001.8
$ws.Range("M107").Value = 1535.66666
$ws.Range("N107").Value = -4841.8
$ws.Range("H132").Value = 2874.4075
$ws.Range("I132").Value = 2026.0588
$ws.Range("K132").Value = 6078.1764
$ws.Range("M132").Value = -3548.1764

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8237.6875
$ws.Range("I7").Value = 10972.429
$ws.Range("J7").Value = 6110.6665
$ws.Range("K7").Value = 10972.429
$ws.Range("L7").Value = 6110.6665
$ws.Range("M7").Value = -10860.429
$ws.Range("N7").Value = -6334.6665
$ws.Range("H16").Value = 377.52942
$ws.Range("I16").Value = 339.46667
$ws.Range("K16").Value = 339.46667
$ws.Range("M16").Value = -169.46667
$ws.Range("H22").Value = 1245.2
$ws.Range("I22").Value = 1182
$ws.Range("K22").Value = 1182
$ws.Range("M22").Value = -887
$ws.Range("H27").Value = 1245.2
$ws.Range("I27").Value = 1182
$ws.Range("K27").Value = 1182
$ws.Range("M27").Value = -1075
$ws.Range("H36").Value = 73238.336
$ws.Range("J36").Value = 73238.336
$ws.Range("L36").Value = 73238.336
$ws.Range("N36").Value = -74362.336
$ws.Range("H40").Value = 14128.077
$ws.Range("I40").Value = 22683
$ws.Range("K40").Value = 22683
$ws.Range("M40").Value = -22547
$ws.Range("H46").Value = 1108.2609
$ws.Range("I46").Value = 415.83334
$ws.Range("K46").Value = 415.83334
$ws.Range("M46").Value = -227.83334
$ws.Range("H82").Value = 2271.6
$ws.Range("I82").Value = 2151.8462
$ws.Range("K82").Value = 2151.8462
$ws.Range("M82").Value = -1790.8462
$ws.Range("H85").Value = 2271.6
$ws.Range("I85").Value = 2151.8462
$ws.Range("K85").Value = 2151.8462
$ws.Range("M85").Value = -903.8462
$ws.Range("H122").Value = 4487.769
$ws.Range("I122").Value = 2450.6
$ws.Range("K122").Value = 7351.799999999999
$ws.Range("M122").Value = -4901.799999999999
$ws.Range("H126").Value = 8237.6875
$ws.Range("I126").Value = 10972.429
$ws.Range("J126").Value = 6110.6665
$ws.Range("K126").Value = 32917.287
$ws.Range("L126").Value = 18331.9995
$ws.Range("M126").Value = -30447.287
$ws.Range("N126").Value = -23271.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 16061.556
$ws.Range("J41").Value = 16017.857
$ws.Range("L41").Value = 16017.857
$ws.Range("N41").Value = -16797.857
$ws.Range("H62").Value = 300
$ws.Range("J62").Value = 300
$ws.Range("L62").Value = 300
$ws.Range("N62").Value = -1548
$ws.Range("H64").Value = 59959
$ws.Range("J64").Value = 59959
$ws.Range("L64").Value = 59959
$ws.Range("N64").Value = -60455
$ws.Range("H65").Value = 300
$ws.Range("J65").Value = 300
$ws.Range("L65").Value = 1500
$ws.Range("N65").Value = -7740
$ws.Range("H67").Value = 59959
$ws.Range("J67").Value = 59959
$ws.Range("L67").Value = 59959
$ws.Range("N67").Value = -61675
$ws.Range("H81").Value = 9595.923000000001
$ws.Range("I81").Value = 35199
$ws.Range("J81").Value = 1915
$ws.Range("K81").Value = 70398
$ws.Range("L81").Value = 3830
$ws.Range("M81").Value = -69337
$ws.Range("N81").Value = -5952
$ws.Range("H84").Value = 9595.923000000001
$ws.Range("I84").Value = 35199
$ws.Range("J84").Value = 1915
$ws.Range("K84").Value = 351990
$ws.Range("L84").Value = 19150
$ws.Range("M84").Value = -346686
$ws.Range("N84").Value = -29758
$ws.Range("H100").Value = 882
$ws.Range("I100").Value = 790.2308
$ws.Range("K100").Value = 1580.4616
$ws.Range("M100").Value = -1039.4616
$ws.Range("H109").Value = 30000
$ws.Range("J109").Value = 30000
$ws.Range("L109").Value = 30000
$ws.Range("N109").Value = -32774
$ws.Range("H122").Value = 5633.25
$ws.Range("I122").Value = 4657
$ws.Range("K122").Value = 13971
$ws.Range("M122").Value = -11521
$ws.Range("H126").Value = 1350.5807
$ws.Range("I126").Value = 1259.8214
$ws.Range("K126").Value = 3779.4642
$ws.Range("M126").Value = -1309.4642
$ws.Range("H132").Value = 2106.5
$ws.Range("I132").Value = 1436.88
$ws.Range("K132").Value = 4310.64
$ws.Range("M132").Value = -1780.64
$ws.Range("H136").Value = 4233.65
$ws.Range("I136").Value = 2670.2
$ws.Range("K136").Value = 8010.599999999999
$ws.Range("M136").Value = -5460.599999999999
